# Update code for mobile domain
# - Adds expected-response sample JSON fragments to the "Get History Completed"
#   and "Get History in Progress" sheets (column D, rows 3-5).
# - Moves the active tab / selection over to the "Get History in Progress" sheet.

$wb = $excel.ActiveWorkbook

# --- "Get History Completed" sheet -----------------------------------------
# NOTE: populate this sheet's new shared strings BEFORE "Get History in
# Progress" so the shared-string table ends up in the same order as the
# authored workbook (COMPLETED entries first, then WAITING entries).
$wsCompleted = $wb.Worksheets.Item("Get History Completed")
$wsCompleted.Range("D3").Value = '"phoneNumber":"081252930398","price":15000,"voucher":0,"status":"COMPLETED"'
$wsCompleted.Range("D4").Value = '"phoneNumber":"0812521617910","price":15000,"voucher":0,"status":"COMPLETED"'
$wsCompleted.Range("D5").Value = '"phoneNumber":"081252161790","price":15000,"voucher":0,"status":"COMPLETED"'
$wsCompleted.Columns.Item(4).ColumnWidth = 103.91673125003601
$wsCompleted.Range("D6").Select()

# --- "Get History in Progress" sheet ----------------------------------------
$wsProgress = $wb.Worksheets.Item("Get History in Progress")
$wsProgress.Range("D3").Value = '"phoneNumber":"081252930398","price":15000,"voucher":0,"status":"WAITING"'
$wsProgress.Range("D4").Value = '"phoneNumber":"0812521617910","price":15000,"voucher":0,"status":"WAITING"'
$wsProgress.Range("D5").Value = '"phoneNumber":"081252161790","price":15000,"voucher":0,"status":"WAITING"'
$wsProgress.Columns.Item(4).ColumnWidth = 68.25006250003605

# Make "Get History in Progress" the active/visible tab with D6 selected.
$wsProgress.Activate()
$wsProgress.Range("D6").Select()
